$wb = $excel.ActiveWorkbook

# xlPasteFormats constant (used below to copy a cell's style without
# touching any cell's stored value).
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. "总计" (summary) sheet: a new 2022-Q3 row is inserted at the top of
#    the data, and every existing data row slides down by one. Easiest
#    to reproduce by rewriting bottom-up so we never clobber a value we
#    still need to read.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("A9").Value = 7
$summary.Range("B9").Value = "2020-Q4"
$summary.Range("C9").Value = 6
$summary.Range("D9").Value = 0.39

$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 6
$summary.Range("D8").Value = 1.09

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 11
$summary.Range("D7").Value = 1.92

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 5
$summary.Range("D6").Value = 1.03

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 4
$summary.Range("D5").Value = 0.78

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 0.29

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.15

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.3

# Row 9 is brand new, so it doesn't inherit the "index column" style (s=2)
# that every other row in column A already carries. Copy it over from the
# row above so A9 matches A2:A8.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q3" worksheet right before "2022-Q2" and
#    fill it with the quarterly fund-holding detail table.
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Header row (B1:H1) + the "index" cell A2 use the same bold/boxed style
# as every other quarter sheet - grab it from the "2022-Q2" sheet instead
# of hand-building a new style entry.
$styleSource = $beforeSheet.Range("B1")

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "000586"
$q3.Range("C2").Value = "景顺长城中小创精选股票"
$q3.Range("D2").Value = "2.21"
$q3.Range("E2").Value = "93.50"
$q3.Range("F2").Value = "7.86"
$q3.Range("G2").Value = "0.1737"
$q3.Range("H2").Value = 6

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "015784"
$q3.Range("C3").Value = "中信建投中证1000指数增强A"
$q3.Range("D3").Value = "8.10"
$q3.Range("E3").Value = "92.20"
$q3.Range("F3").Value = "0.65"
$q3.Range("G3").Value = "0.0526"
$q3.Range("H3").Value = 8

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "260115"
$q3.Range("C4").Value = "景顺长城中小盘混合"
$q3.Range("D4").Value = "0.92"
$q3.Range("E4").Value = "92.87"
$q3.Range("F4").Value = "5.27"
$q3.Range("G4").Value = "0.0485"
$q3.Range("H4").Value = 7

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "015785"
$q3.Range("C5").Value = "中信建投中证1000指数增强C"
$q3.Range("D5").Value = "3.32"
$q3.Range("E5").Value = "92.20"
$q3.Range("F5").Value = "0.65"
$q3.Range("G5").Value = "0.0216"
$q3.Range("H5").Value = 8

# The fund-code / size / position columns look numeric ("000586", "2.21",
# "93.50", ...) so a plain .Value assignment gets auto-coerced to a
# number (and drops the leading zero on the fund code). Force text entry
# the same way Excel's "Format Cells > Text" + retype does, then strip
# the temporary number-format back off via a formats-only paste from a
# plain cell so the result matches a plain inlineStr cell with no style.
$textValues = @{
    "B2" = "000586";     "C2" = "景顺长城中小创精选股票"; "D2" = "2.21";  "E2" = "93.50"; "F2" = "7.86"; "G2" = "0.1737"
    "B3" = "015784";     "C3" = "中信建投中证1000指数增强A"; "D3" = "8.10";  "E3" = "92.20"; "F3" = "0.65"; "G3" = "0.0526"
    "B4" = "260115";     "C4" = "景顺长城中小盘混合"; "D4" = "0.92";  "E4" = "92.87"; "F4" = "5.27"; "G4" = "0.0485"
    "B5" = "015785";     "C5" = "中信建投中证1000指数增强C"; "D5" = "3.32";  "E5" = "92.20"; "F5" = "0.65"; "G5" = "0.0216"
}
foreach ($addr in $textValues.Keys) {
    $cell = $q3.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$addr]
}
$blank = $summary.Range("D2")
$blank.Copy()
$q3.Range("B2:G5").PasteSpecial($xlPasteFormats)

# Re-apply the real header/index style now that the formats-only paste
# above reset B1:H1 / A2 back to the default style.
$styleSource.Copy()
$q3.Range("B1:H1").PasteSpecial($xlPasteFormats)
$summary.Range("A2").Copy()
$q3.Range("A2").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 3. Restore the original "active sheet" (last tab, 2020-Q4) so the
#    freshly-inserted sheet doesn't steal tab-selection / activeTab.
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()
